$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export Worksheet")

# Update the header for column G
$ws.Range("G1").Value = "End PM"

# End PM values keyed by row number
$endPm = @{
    2 = 33.0
    3 = 15.0
    4 = 34.8
    5 = 16.2
    6 = 17.4
    7 = 26.3
    8 = 10.726
    9 = 36.4
    10 = 19.0
    11 = 21.6
    12 = 13.579
    13 = 19.0
    14 = 35.1
    15 = 42.014
    16 = 28.923
    17 = 48.558
    18 = 10.053
    19 = 37.9
}

foreach ($row in $endPm.Keys) {
    $val = $endPm[$row]
    $text = "End PM: {0,8:F3}" -f $val
    $cell = $ws.Cells.Item($row, 7)
    $cell.ClearFormats()
    $cell.Value = $text
}

# Update the SQL text in the SQL worksheet to also select/alias the end_pm
# column the same way the beg_pm column is already formatted.
$sqlWs = $wb.Worksheets.Item("SQL")
$newSql = "select a.ea, a.treatment, a.county, a.route, a.year, ('Beg PM: ' || to_char(a.beg_pm, 990.999)) as ""Beg PM"", ('End PM: ' || to_char(a.end_pm, 990.999)) as ""End PM"", (a.end_pm-a.beg_pm) as length, a.budget_group from s1383currentr a 
where a.county = 'SM'  
union  
select b.ea, b.treatment, b.county, b.route, b.year, ('Beg PM: ' || to_char(b.beg_pm, 990.999)) as ""Beg PM"",  ('End PM: ' || to_char(b.end_pm, 990.999)) as ""End PM"",  (b.end_pm-b.beg_pm) as length, b.budget_group from s1383historyr b 
where b.county = 'SM'  
order by year"
$sqlWs.Range("A2").Value = $newSql
